$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/6/2025  Through  1/12/2025"

# --- Data cell updates ---
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("D15").Copy($ws.Range("J15"))
$ws.Range("J15").Value = 1
$ws.Range("H14").Copy($ws.Range("K15"))
$ws.Range("K15").Value = -100
$ws.Range("C16").Value = 4
$ws.Range("D15").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 4
$ws.Range("H14").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -11.111111111111
$ws.Range("I16").Value = 5
$ws.Range("D15").Copy($ws.Range("J16"))
$ws.Range("J16").Value = 4
$ws.Range("H14").Copy($ws.Range("K16"))
$ws.Range("K16").Value = 25
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -72.222222222222
$ws.Range("N16").Value = -86.842105263157
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = -38.461538461538
$ws.Range("L17").Value = -11.111111111111
$ws.Range("N17").Value = -66.666666666666
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("D15").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1
$ws.Range("H14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -40
$ws.Range("D15").Copy($ws.Range("J18"))
$ws.Range("J18").Value = 1
$ws.Range("H14").Copy($ws.Range("K18"))
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = -92.307692307692
$ws.Range("N18").Value = -97.5
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 80
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 50
$ws.Range("I19").Value = 11
$ws.Range("J19").Value = 7
$ws.Range("K19").Value = 57.142857142857
$ws.Range("L19").Value = -15.384615384615
$ws.Range("M19").Value = -15.384615384615
$ws.Range("N19").Value = -83.582089552238
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = 57.142857142857
$ws.Range("M20").Value = -21.428571428571
$ws.Range("N20").Value = -82.539682539682
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 36.842105263157
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 86
$ws.Range("H21").Value = -10.465116279069
$ws.Range("I21").Value = 36
$ws.Range("J21").Value = 36
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 9.090909090909
$ws.Range("M21").Value = -45.454545454545
$ws.Range("N21").Value = -85.123966942148
$ws.Range("D15").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("H14").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -100
$ws.Range("D15").Copy($ws.Range("G23"))
$ws.Range("G23").Value = 1
$ws.Range("H14").Copy($ws.Range("H23"))
$ws.Range("H23").Value = -100
$ws.Range("D15").Copy($ws.Range("J23"))
$ws.Range("J23").Value = 1
$ws.Range("H14").Copy($ws.Range("K23"))
$ws.Range("K23").Value = -100
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -51.851851851851
$ws.Range("F24").Value = 70
$ws.Range("H24").Value = -21.348314606741
$ws.Range("I24").Value = 19
$ws.Range("J24").Value = 42
$ws.Range("K24").Value = -54.761904761904
$ws.Range("L24").Value = -40.625
$ws.Range("M24").Value = -47.222222222222
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -83.333333333333
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -63.888888888888
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 23
$ws.Range("K25").Value = -78.260869565217
$ws.Range("L25").Value = -50
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 57
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = 23.913043478260
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = 62.5
$ws.Range("L26").Value = 52.941176470588
$ws.Range("M26").Value = 23.809523809523
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("D15").Copy($ws.Range("J27"))
$ws.Range("J27").Value = 1
$ws.Range("H14").Copy($ws.Range("K27"))
$ws.Range("K27").Value = -100
$ws.Range("D15").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 3
$ws.Range("H14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -66.666666666666
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("D15").Copy($ws.Range("I28"))
$ws.Range("I28").Value = 1
$ws.Range("D15").Copy($ws.Range("J28"))
$ws.Range("J28").Value = 3
$ws.Range("H14").Copy($ws.Range("K28"))
$ws.Range("K28").Value = -66.666666666666
$ws.Range("H14").Copy($ws.Range("L28"))
$ws.Range("L28").Value = 0
$ws.Range("G29").Value = 5
$ws.Range("G30").Value = 2
$ws.Range("H14").Copy($ws.Range("L33"))
$ws.Range("L33").Value = -100
$ws.Range("J42").Value = 353
$ws.Range("K42").Value = -11.970074812967
$ws.Range("L42").Value = -27.663934426229
$ws.Range("M42").Value = -53.430079155672
$ws.Range("N42").Value = -53.552631578947
$ws.Range("J43").Value = 73
$ws.Range("K43").Value = -84.164859002169
$ws.Range("L43").Value = -88.786482334869
$ws.Range("M43").Value = -93.652173913043
$ws.Range("N43").Value = -95.400126023944
$ws.Range("J46").Value = 1072
$ws.Range("K46").Value = -59.547169811320
$ws.Range("L46").Value = -70.378557612600
$ws.Range("M46").Value = -88.111345236775
$ws.Range("N46").Value = -88.373101952277
